$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 44
$ws.Range("H44").Value = 44025
$ws.Range("J44").Value = 44025
$ws.Range("L44").Value = 44025
$ws.Range("N44").Value = -44949

# Row 86
$ws.Range("H86").Value = 112563.91
$ws.Range("I86").Value = 205167.17
$ws.Range("J86").Value = 1440
$ws.Range("K86").Value = 205167.17
$ws.Range("L86").Value = 1440
$ws.Range("M86").Value = -204044.17
$ws.Range("N86").Value = -3686

# Row 89
$ws.Range("H89").Value = 112563.91
$ws.Range("I89").Value = 205167.17
$ws.Range("J89").Value = 1440
$ws.Range("K89").Value = 1025835.85
$ws.Range("L89").Value = 7200
$ws.Range("M89").Value = -1020219.85
$ws.Range("N89").Value = -18432

# Row 106
$ws.Range("H106").Value = 2466.5715
$ws.Range("I106").Value = 1682.3636
$ws.Range("J106").Value = 3329.2
$ws.Range("K106").Value = 1682.3636
$ws.Range("L106").Value = 3329.2
$ws.Range("M106").Value = -1051.3636
$ws.Range("N106").Value = -4591.2

# Row 125
$ws.Range("H125").Value = 38209.68
$ws.Range("I125").Value = 200772.8
$ws.Range("J125").Value = 2869.8696
$ws.Range("K125").Value = 1806955.2
$ws.Range("L125").Value = 25828.8264
$ws.Range("M125").Value = -1804495.2
$ws.Range("N125").Value = -30748.8264

# Row 127
$ws.Range("H127").Value = 461.86206
$ws.Range("I127").Value = 245.58333
$ws.Range("J127").Value = 1500
$ws.Range("K127").Value = 736.74999
$ws.Range("L127").Value = 4500
$ws.Range("M127").Value = 4223.25001
$ws.Range("N127").Value = -14420

# Row 135
$ws.Range("H135").Value = 928.5
$ws.Range("I135").Value = 487.2
$ws.Range("J135").Value = 1329.6818
$ws.Range("K135").Value = 4384.8
$ws.Range("L135").Value = 11967.1362
$ws.Range("M135").Value = -1849.8
$ws.Range("N135").Value = -17037.1362

# Row 138
$ws.Range("H138").Value = 1252.3043
$ws.Range("I138").Value = 1170.0333
$ws.Range("J138").Value = 1406.5625
$ws.Range("K138").Value = 3510.0999
$ws.Range("L138").Value = 4219.6875
$ws.Range("M138").Value = 1629.9001
$ws.Range("N138").Value = -14499.6875

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 16699.033
$ws.Range("I32").Value = 9896.115
$ws.Range("K32").Value = 9896.115
$ws.Range("M32").Value = -9609.115

# Row 45
$ws.Range("H45").Value = 1370.7142
$ws.Range("I45").Value = 1243
$ws.Range("J45").Value = 1498.4286
$ws.Range("K45").Value = 1243
$ws.Range("L45").Value = 1498.4286
$ws.Range("M45").Value = -866
$ws.Range("N45").Value = -2252.4286

# Row 132
$ws.Range("H132").Value = 1332604
$ws.Range("I132").Value = 1880.5435
$ws.Range("J132").Value = 2607880.8
$ws.Range("K132").Value = 5641.6305
$ws.Range("L132").Value = 7823642.399999999
$ws.Range("M132").Value = -3111.6305
$ws.Range("N132").Value = -7828702.399999999

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1007.75
$ws.Range("I86").Value = 991.75
$ws.Range("J86").Value = 1071.75
$ws.Range("K86").Value = 991.75
$ws.Range("L86").Value = 1071.75
$ws.Range("M86").Value = 131.25
$ws.Range("N86").Value = -3317.75

# Row 89
$ws.Range("H89").Value = 1007.75
$ws.Range("I89").Value = 991.75
$ws.Range("J89").Value = 1071.75
$ws.Range("K89").Value = 4958.75
$ws.Range("L89").Value = 5358.75
$ws.Range("M89").Value = 657.25
$ws.Range("N89").Value = -16590.75

# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# Row 134
$ws.Range("H134").Value = 4283.5425
$ws.Range("I134").Value = 1820.0312
$ws.Range("J134").Value = 7203.2593
$ws.Range("K134").Value = 5460.0936
$ws.Range("L134").Value = 21609.7779
$ws.Range("M134").Value = -2925.0936
$ws.Range("N134").Value = -26679.7779

$ws = $wb.Worksheets.Item("CRP")
# Row 103
$ws.Range("H103").Value = 12408.223
$ws.Range("I103").Value = 2381
$ws.Range("J103").Value = 20430
$ws.Range("K103").Value = 2381
$ws.Range("L103").Value = 20430
$ws.Range("M103").Value = -1209
$ws.Range("N103").Value = -22774

# Row 122
$ws.Range("H122").Value = 1681.8064
$ws.Range("I122").Value = 1046.4445
$ws.Range("J122").Value = 2561.5386
$ws.Range("K122").Value = 3139.3335
$ws.Range("L122").Value = 7684.6158
$ws.Range("M122").Value = -689.3335000000002
$ws.Range("N122").Value = -12584.6158

$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 2136
$ws.Range("J80").Value = 2626.6667
$ws.Range("L80").Value = 7880.000100000001
$ws.Range("N80").Value = -9752.000100000001

# Row 83
$ws.Range("H83").Value = 2136
$ws.Range("J83").Value = 2626.6667
$ws.Range("L83").Value = 23640.0003
$ws.Range("N83").Value = -33000.0003

# Row 115
$ws.Range("H115").Value = 1345.3334
$ws.Range("I115").Value = 1010.2727
$ws.Range("J115").Value = 5031
$ws.Range("K115").Value = 3030.8181
$ws.Range("L115").Value = 15093
$ws.Range("M115").Value = -1855.8181
$ws.Range("N115").Value = -17443

# Row 134
$ws.Range("H134").Value = 3408.4333
$ws.Range("I134").Value = 1806.6666
$ws.Range("J134").Value = 4094.9048
$ws.Range("K134").Value = 5419.9998
$ws.Range("L134").Value = 12284.7144
$ws.Range("M134").Value = -349.9997999999996
$ws.Range("N134").Value = -22424.7144

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2249.6924
$ws.Range("I7").Value = 2248.889
$ws.Range("J7").Value = 2251.5
$ws.Range("K7").Value = 2248.889
$ws.Range("L7").Value = 2251.5
$ws.Range("M7").Value = -2136.889
$ws.Range("N7").Value = -2475.5

# Row 40
$ws.Range("H40").Value = 2346.325
$ws.Range("I40").Value = 1935.6333
$ws.Range("K40").Value = 1935.6333
$ws.Range("M40").Value = -1799.6333

# Row 46
$ws.Range("H46").Value = 1668012.1
$ws.Range("I46").Value = 1740.5
$ws.Range("J46").Value = 2501148
$ws.Range("K46").Value = 1740.5
$ws.Range("L46").Value = 2501148
$ws.Range("M46").Value = -1552.5
$ws.Range("N46").Value = -2501524

# Row 61
$ws.Range("H61").Value = 2250.1516
$ws.Range("I61").Value = 2341.1304
$ws.Range("J61").Value = 2040.9
$ws.Range("K61").Value = 2341.1304
$ws.Range("L61").Value = 2040.9
$ws.Range("M61").Value = -2139.1304
$ws.Range("N61").Value = -2444.9

# Row 100
$ws.Range("H100").Value = 2268.0588
$ws.Range("I100").Value = 2000.909
$ws.Range("J100").Value = 2757.8333
$ws.Range("K100").Value = 2000.909
$ws.Range("L100").Value = 2757.8333
$ws.Range("M100").Value = -1459.909
$ws.Range("N100").Value = -3839.8333

# Row 113
$ws.Range("H113").Value = 2250.1516
$ws.Range("I113").Value = 2341.1304
$ws.Range("J113").Value = 2040.9
$ws.Range("K113").Value = 2341.1304
$ws.Range("L113").Value = 2040.9
$ws.Range("M113").Value = -171.1304
$ws.Range("N113").Value = -6380.9

# Row 126
$ws.Range("H126").Value = 2249.6924
$ws.Range("I126").Value = 2248.889
$ws.Range("J126").Value = 2251.5
$ws.Range("K126").Value = 6746.667
$ws.Range("L126").Value = 6754.5
$ws.Range("M126").Value = -4276.667
$ws.Range("N126").Value = -11694.5

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1463.746
$ws.Range("I132").Value = 1045.6052
$ws.Range("J132").Value = 2099.32
$ws.Range("K132").Value = 3136.8156
$ws.Range("L132").Value = 6297.960000000001
$ws.Range("M132").Value = -606.8155999999999
$ws.Range("N132").Value = -11357.96

# Row 136
$ws.Range("H136").Value = 1041.2439
$ws.Range("I136").Value = 731.9149
$ws.Range("J136").Value = 1456.6285
$ws.Range("K136").Value = 2195.7447
$ws.Range("L136").Value = 4369.8855
$ws.Range("M136").Value = 354.2552999999998
$ws.Range("N136").Value = -9469.8855
